# Updates the cryptos list values (Price + Volume(1h)) to the latest snapshot,
# matching the "Updated cryptos list ... with GitHub Actions" refresh job.
# Rows 48/49 swap the ARBITRUM / Filecoin entries (Filecoin now ranks above
# ARBITRUM) and row 51 replaces Optimism with Cronos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is an A1 cell + its new value. Column D (Price) values that look
# like a plain number (e.g. "585.45") are written with a leading apostrophe so
# Excel keeps them as text -- exactly like the source data -- instead of silently
# re-typing the cell as a Number.
$updates = @(
    @{ Cell = "D2"; Value = "67.253.43" }
    @{ Cell = "E2"; Value = "  +0.66%  " }
    @{ Cell = "D3"; Value = "2.493.60" }
    @{ Cell = "E3"; Value = "  +0.68%  " }
    @{ Cell = "E4"; Value = "  +0.00%  " }
    @{ Cell = "D5"; Value = "'585.45" }
    @{ Cell = "E5"; Value = "  +0.32%  " }
    @{ Cell = "D6"; Value = "'172.37" }
    @{ Cell = "E6"; Value = "  +3.10%  " }
    @{ Cell = "E7"; Value = "  -0.12%  " }
    @{ Cell = "D8"; Value = "'0.514" }
    @{ Cell = "E8"; Value = "  -0.33%  " }
    @{ Cell = "D9"; Value = "2.494.06" }
    @{ Cell = "E9"; Value = "  +0.72%  " }
    @{ Cell = "D10"; Value = "'0.136" }
    @{ Cell = "E10"; Value = "  +0.79%  " }
    @{ Cell = "E11"; Value = "  +0.19%  " }
    @{ Cell = "E12"; Value = "  +0.15%  " }
    @{ Cell = "E13"; Value = "  -0.71%  " }
    @{ Cell = "D14"; Value = "'25.52" }
    @{ Cell = "E14"; Value = "  -1.14%  " }
    @{ Cell = "D15"; Value = "2.919.63" }
    @{ Cell = "D16"; Value = "67.176.34" }
    @{ Cell = "E16"; Value = "  +0.63%  " }
    @{ Cell = "E17"; Value = "  -1.46%  " }
    @{ Cell = "D18"; Value = "2.498.05" }
    @{ Cell = "E18"; Value = "  +0.89%  " }
    @{ Cell = "D19"; Value = "'11.06" }
    @{ Cell = "E19"; Value = "  -4.66%  " }
    @{ Cell = "D20"; Value = "'7.46" }
    @{ Cell = "E20"; Value = "  -5.06%  " }
    @{ Cell = "D21"; Value = "'351.44" }
    @{ Cell = "E21"; Value = "  -2.89%  " }
    @{ Cell = "D22"; Value = "'4.04" }
    @{ Cell = "E22"; Value = "  -0.28%  " }
    @{ Cell = "D23"; Value = "'1.00" }
    @{ Cell = "E23"; Value = "  -0.04%  " }
    @{ Cell = "E24"; Value = "  -4.36%  " }
    @{ Cell = "D25"; Value = "'68.68" }
    @{ Cell = "E25"; Value = "  -3.03%  " }
    @{ Cell = "E26"; Value = "  -1.87%  " }
    @{ Cell = "E27"; Value = "  -1.89%  " }
    @{ Cell = "D28"; Value = "'1.00" }
    @{ Cell = "E28"; Value = "  +0.30%  " }
    @{ Cell = "D29"; Value = "2.622.61" }
    @{ Cell = "E29"; Value = "  +0.75%  " }
    @{ Cell = "D30"; Value = "0.0₃0904" }
    @{ Cell = "E30"; Value = "  -2.80%  " }
    @{ Cell = "D31"; Value = "'512.31" }
    @{ Cell = "E31"; Value = "  -0.08%  " }
    @{ Cell = "D32"; Value = "'7.84" }
    @{ Cell = "E32"; Value = "  -2.69%  " }
    @{ Cell = "E33"; Value = "  -2.10%  " }
    @{ Cell = "E34"; Value = "  -3.01%  " }
    @{ Cell = "E35"; Value = "  -0.03%  " }
    @{ Cell = "D36"; Value = "'160.26" }
    @{ Cell = "E36"; Value = "  +1.20%  " }
    @{ Cell = "E37"; Value = "  -6.78%  " }
    @{ Cell = "E38"; Value = "  +0.85%  " }
    @{ Cell = "D39"; Value = "'18.28" }
    @{ Cell = "E39"; Value = "  -3.27%  " }
    @{ Cell = "E40"; Value = "  -5.37%  " }
    @{ Cell = "E41"; Value = "  -2.93%  " }
    @{ Cell = "E42"; Value = "  -0.12%  " }
    @{ Cell = "E43"; Value = "  -1.80%  " }
    @{ Cell = "E44"; Value = "  -1.20%  " }
    @{ Cell = "D45"; Value = "'2.37" }
    @{ Cell = "E45"; Value = "  -3.04%  " }
    @{ Cell = "D46"; Value = "'38.81" }
    @{ Cell = "D47"; Value = "'143.40" }
    @{ Cell = "E47"; Value = "  +0.79%  " }
    @{ Cell = "B48"; Value = "Filecoin" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" }
    @{ Cell = "D48"; Value = "'3.46" }
    @{ Cell = "E48"; Value = "  -3.48%  " }
    @{ Cell = "B49"; Value = "ARBITRUM" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" }
    @{ Cell = "D49"; Value = "'0.516" }
    @{ Cell = "E49"; Value = "  -4.09%  " }
    @{ Cell = "D50"; Value = "0.0₆0252" }
    @{ Cell = "E50"; Value = "  -5.72%  " }
    @{ Cell = "B51"; Value = "Cronos" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" }
    @{ Cell = "D51"; Value = "'0.0731" }
    @{ Cell = "E51"; Value = "  -0.83%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
